$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"

# Row 3
$ws.Range("B3").Value = "[-, 'MEC-3B-Calderaria', -, -]"
$ws.Range("C3").Value = "[-, -, -, 'MCT-1A-Metrologia']"
$ws.Range("D3").Value = "-"

# Row 4
$ws.Range("B4").Value = "[-, 'MEC-3B-Calderaria', -, -]"
$ws.Range("C4").Value = "[-, -, -, 'MCT-1A-Metrologia']"
$ws.Range("D4").Value = "-"

# Row 6
$ws.Range("B6").Value = "[-, 'MEC-3B-Calderaria', -, -]"
$ws.Range("C6").Value = "[-, -, -, 'MCT-1A-Metrologia']"
$ws.Range("D6").Value = "-"

# Row 7
$ws.Range("B7").Value = "[-, 'MEC-3B-Calderaria', -, -]"
$ws.Range("C7").Value = "[-, -, -, 'MCT-1A-Metrologia']"

# Row 11
$ws.Range("C11").Value = "-"
$ws.Range("E11").Value = "[-, -, -, 'MEC-3A-Calderaria']"

# Row 12
$ws.Range("C12").Value = "-"
$ws.Range("E12").Value = "[-, -, -, 'MEC-3A-Calderaria']"

# Row 14
$ws.Range("C14").Value = "-"
$ws.Range("E14").Value = "[-, -, -, 'MEC-3A-Calderaria']"

# Row 15
$ws.Range("C15").Value = "-"
$ws.Range("E15").Value = "[-, -, -, 'MEC-3A-Calderaria']"

# Row 18
$ws.Range("C18").Value = "[-, -, 'MEC-1NA-Tec. Mat. Não Metal.', -]"
$ws.Range("D18").Value = "[-, -, -, 'MEC-1NB-Caldeiraria']"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "[-, -, -, 'MEC-1NB-Caldeiraria']"

# Row 19
$ws.Range("C19").Value = "[-, -, 'MEC-1NA-Tec. Mat. Não Metal.', -]"
$ws.Range("F19").Value = "[-, -, -, 'MEC-1NB-Caldeiraria']"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "[-, -, 'MEC-1NA-Tec. Mat. Não Metal.', -]"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "[-, -, 'MEC-1NA-Tec. Mat. Não Metal.', -]"
$ws.Range("F21").Value = "[-, -, -, 'MEC-1NB-Caldeiraria']"
